# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# The underlying data rows for a handful of matches were recorded in the
# wrong order (e.g. the "home" and "away" fixture rows were transposed).
# This script swaps the full data payload (columns B through AC) between
# each pair of rows, leaving column A (the running row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of row numbers whose contents (columns B:AC) must be swapped.
$rowPairs = @(
    @(213, 214),
    @(215, 216),
    @(229, 231),
    @(232, 233),
    @(251, 252),
    @(263, 265)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AC$r1")
    $range2 = $ws.Range("B$r2" + ":AC$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
